$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5 (existing rows 5.. shift down by one)
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new record
$ws.Cells.Item(5, 1).Value = 7
$ws.Cells.Item(5, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(5, 3).Value = "Ñuble"
$ws.Cells.Item(5, 4).Value = Get-Date -Year 2023 -Month 6 -Day 6 -Hour 0 -Minute 0 -Second 0
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 100112013
$ws.Cells.Item(5, 7).Value = "Alcachofa"
$ws.Cells.Item(5, 8).Value = "Española"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 20
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 16000
$ws.Cells.Item(5, 13).Value = 16000
$ws.Cells.Item(5, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(5, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(5, 16).Value = 267
$ws.Cells.Item(5, 17).Value = 60
$ws.Cells.Item(5, 18).Value = "Hortaliza"
